$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update claim numbers (NroSiniestro) in column F (leading apostrophe keeps them as text,
# preserving the existing "text" cell style/quote-prefix formatting)
$ws.Range("F2").Value = "'1220194200667"
$ws.Range("F3").Value = "'1120194100412"
$ws.Range("F4").Value = "'0420194406717"

# Update user (Usuario) in column C for row 4
$ws.Range("C4").Value = "apellegrini"

# Update the active selection to F4
$ws.Range("F4").Select()
